$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 3
    4  = 1
    7  = 5
    9  = -2
    10 = -6
    11 = 2
    12 = 3
    13 = 1
    14 = -3
    16 = -2
    17 = -6
    18 = -1
    19 = -2
    22 = 1
    24 = -1
    25 = 10
    26 = 1
    27 = -7
    28 = -4
    29 = 3
    30 = 11
    31 = -1
    32 = 1
    33 = 2
    34 = -1
    36 = -2
    38 = -7
    39 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
